# [Kadastro App] Yeni kayit eklendi: 3012
# Adds the new record (Kayit No 3012) as row 71 on both the master
# "Kayitlar" sheet and the filtered "Erdemli" district sheet.

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "3012"
    B = "2025-09-11"
    C = "Erdemli"
    D = "1"
    E = "3B"
    F = "ÖZKAN AKBAŞ (Mühendis), SERDAR ARSLAN (Tekniker)"
}

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Keep the appended row text-formatted (matches every other row in
    # the sheet, where "numeric-looking" values such as Kayit No / Parsel
    # Sayisi are stored as text).
    $rowRange = $ws.Range("A71:F71")
    $rowRange.NumberFormat = "@"

    $ws.Range("A71").Value = $newRow.A
    $ws.Range("B71").Value = $newRow.B
    $ws.Range("C71").Value = $newRow.C
    $ws.Range("D71").Value = $newRow.D
    $ws.Range("E71").Value = $newRow.E
    $ws.Range("F71").Value = $newRow.F
}
